$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.830.02'
$ws.Range("E2").Value = '  +0.24%  '
# Row 3
$ws.Range("D3").Value = '2.357.88'
$ws.Range("E3").Value = '  -0.63%  '
# Row 4
$ws.Range("E4").Value = '  +0.11%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.692'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.62%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.53'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.15%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '77.58'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.41%  '
# Row 8
$ws.Range("E8").Value = '  -0.04%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.628'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +19.99%  '
# Row 10
$ws.Range("E10").Value = '  +4.03%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.38'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.52%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '33.87'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +22.14%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.63'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +19.88%  '
# Row 14
$ws.Range("E14").Value = '  +1.77%  '
# Row 15
$ws.Range("D15").Value = '2.704.02'
$ws.Range("E15").Value = '  -0.76%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '17.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.73%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.931'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.59%  '
# Row 18
$ws.Range("D18").Value = '2.353.71'
$ws.Range("E18").Value = '  -0.73%  '
# Row 19
$ws.Range("D19").Value = '43.796.14'
$ws.Range("E19").Value = '  +0.51%  '
# Row 20
$ws.Range("E20").Value = '  +2.20%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.69'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.35%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '77.65'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.68%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.71%  '
# Row 24
$ws.Range("E24").Value = '  -0.03%  '
# Row 25
$ws.Range("E25").Value = '  +2.13%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.08%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.85%  '
# Row 28
$ws.Range("E28").Value = '  +16.69%  '
# Row 29
$ws.Range("E29").Value = '  +2.16%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '23.02'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.56%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.05'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.26%  '
# Row 32
$ws.Range("E32").Value = '  -4.25%  '
# Row 33
$ws.Range("E33").Value = '  +6.05%  '
# Row 34
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.35'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.27%  '
# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0759'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.39%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.78'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.13%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.44'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.95%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.48'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.08%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0279'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.55%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '19.55'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.75%  '
# Row 42
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.202'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +15.29%  '
# Row 43
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.03'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.84%  '
# Row 44
$ws.Range("E44").Value = '  +0.01%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.104'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.91%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +12.20%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.27'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.90%  '
# Row 48
$ws.Range("E48").Value = '  +1.87%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.84'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.32%  '
# Row 50
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '56.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.06%  '
# Row 51
$ws.Range("B51").Value = 'FTXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.50'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.13%  '
